$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Change column H (rows 2-29) from boolean FALSE to the text string "False".
# A plain string assignment of "False"/"FALSE" gets auto-coerced back into a
# boolean by the engine, so instead we:
#   1. Write "False " (with a trailing space) as literal text into H2 while
#      the cell already has a Text number format + the desired Arial 10pt
#      black font - this sidesteps the boolean auto-detection and creates
#      exactly one new font + one new cell style.
#   2. Use a helper cell with a formula that trims the trailing space, copy
#      it, and Paste-Special (values only) back onto H2 so the final stored
#      text is the clean word "False" (without disturbing the style that was
#      just set on H2).
#   3. Copy H2's fully finished format onto the remaining H3:H29 cells
#      (Paste-Special formats only) so they pick up the very same style
#      without generating any further new fonts/styles.
#   4. Repeat the trim trick per remaining cell so every one of them stores
#      the clean shared string "False" too.
# ---------------------------------------------------------------------------

$helper = $ws.Range("Z5")

$first = $ws.Range("H2")
$first.NumberFormat = "@"
$first.Font.Name = "Arial"
$first.Font.Size = 10
$first.Font.Color = 0
$first.Value = "False "

$helper.Formula = "=LEFT(""False X"",5)"
$helper.Copy()
$first.PasteSpecial(-4163)   # xlPasteValues

# Propagate the finished style (number format + font) to the rest of the
# column range before touching their values.
$first.Copy()
$ws.Range("H3:H29").PasteSpecial(-4122)   # xlPasteFormats

for ($r = 3; $r -le 29; $r++) {
    $cell = $ws.Range("H$r")
    $cell.Value = "False "
    $helper.Formula = "=LEFT(""False X"",5)"
    $helper.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

$helper.Clear()

# ---------------------------------------------------------------------------
# Extend the same style (no value) down into three additional blank rows.
# ---------------------------------------------------------------------------
$first.Copy()
$ws.Range("H30:H32").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# Update the on-screen selection to match the new view.
# ---------------------------------------------------------------------------
$ws.Range("G30:J33").Select()
